# v2.0 Fix push button
# Update the "Where to Buy" link for the OFF-(ON) N/O Push-Button Switch row (row 12)
# to point at the new eBay listing, and refresh the sheet's saved scroll/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "http://www.ebay.com/itm/10-Pcs-AC-125V-3A-250V-1-5A-SPST-On-Off-Latching-Red-Push-Button-Switch-/321195563845?hash=item4ac8bf5b45:g:rToAAOSwn7JYDood"

$cell = $ws.Range("C12")

# Remove the old hyperlink (if any) before adding the new one.
if ($cell.Hyperlinks.Count -gt 0) {
    $cell.Hyperlinks.Delete()
}

$ws.Hyperlinks.Add($cell, $newUrl, [Type]::Missing, [Type]::Missing, $newUrl) | Out-Null

# Row height for the push-button row grew to accommodate the longer URL text.
$ws.Rows.Item(12).RowHeight = 65

# Restore the view to the state captured after the edit (scrolled up a couple rows,
# with C12 as the active/selected cell).
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C12").Select()
